$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PDiBCpDoC")

# --- "PDiBCpDoC" sheet updates ---
# B2: updated percentage decline value
$wsData.Range("B2").Value = 0.17

[void]$wsData.Activate()
[void]$wsData.Range("B3").Select()

# --- "About" sheet updates ---
# B4: reporting year 2019 -> 2024
$wsAbout.Range("B4").Value = 2024

# B5/B6: source title & URL refreshed to the new citation
$wsAbout.Range("B5").Value = "Electric Vehicle Outlook 2024"
$wsAbout.Range("B6").Value = "https://about.bnef.com/electric-vehicle-outlook/"

# Row 14 (D14) placeholder cell is no longer part of the used range
[void]$wsAbout.Range("D14").Clear()

# "About" stays the tab that's active/selected when the workbook is saved
[void]$wsAbout.Activate()
[void]$wsAbout.Range("B6").Select()
